$d = $word.ActiveDocument

$replacements = @(
    @{old="492×7="; new="477×4="},
    @{old="191×9="; new="232×8="},
    @{old="800×4="; new="465×3="},
    @{old="114×5="; new="751×5="},
    @{old="166×4="; new="315×4="},
    @{old="931×3="; new="514×8="},
    @{old="583×2="; new="482×6="},
    @{old="531×2="; new="429×8="},
    @{old="634×6="; new="931×7="},
    @{old="307×5="; new="782×8="},
    @{old="604×5="; new="595×9="},
    @{old="614×9="; new="222×3="},
    @{old="321×9="; new="150×3="},
    @{old="874×5="; new="909×9="},
    @{old="825×6="; new="672×7="},
    @{old="946×7="; new="535×5="},
    @{old="203×4="; new="348×2="},
    @{old="629×9="; new="498×8="},
    @{old="831×8="; new="488×8="},
    @{old="103×3="; new="211×7="},
    @{old="845×8="; new="130×4="},
    @{old="914×2="; new="206×7="},
    @{old="253×6="; new="273×2="},
    @{old="684×9="; new="923×9="},
    @{old="738×4="; new="469×2="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
